$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 value change (A2 stays 7)
$ws.Range("B2").Value = -1.169

# Row 3: A3 13 -> 9, B3 0.265 -> 0.341
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = 0.341

# Row 4: A4 9 -> 13, B4 0.334 -> 0.372
$ws.Range("A4").Value = 13
$ws.Range("B4").Value = 0.372

# Row 5: B5 0.697 -> 0.73 (A5 stays 3)
$ws.Range("B5").Value = 0.73

# Row 6: B6 1.361 -> 1.289 (A6 stays 4)
$ws.Range("B6").Value = 1.289
